$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.516.46"
$ws.Range("E2").Value = "  -1.30%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.854.55"
$ws.Range("E3").Value = "  -2.30%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "520.79"
$ws.Range("E5").Value = "  +5.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.35"
$ws.Range("E6").Value = "  -4.48%  "
$ws.Range("E7").Value = "  -2.73%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("E9").Value = "  -3.52%  "
$ws.Range("E10").Value = "  -6.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000322"
$ws.Range("E11").Value = "  -8.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "41.50"
$ws.Range("E12").Value = "  -3.77%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.41"
$ws.Range("E13").Value = "  -0.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.468.63"
$ws.Range("E14").Value = "  -2.35%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.22"
$ws.Range("E15").Value = "  +6.48%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.845.78"
$ws.Range("E16").Value = "  -3.86%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.05"
$ws.Range("E17").Value = "  -1.79%  "
$ws.Range("E18").Value = "  -2.21%  "
$ws.Range("E19").Value = "  +2.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "68.512.32"
$ws.Range("E20").Value = "  -1.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "415.72"
$ws.Range("E21").Value = "  -5.72%  "
$ws.Range("E22").Value = "  +0.53%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.96"
$ws.Range("E23").Value = "  -3.79%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "86.81"
$ws.Range("E24").Value = "  -2.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.97"
$ws.Range("E25").Value = "  +6.12%  "
$ws.Range("E26").Value = "  -2.81%  "
$ws.Range("E27").Value = "  -5.80%  "
$ws.Range("E28").Value = "  -4.74%  "
$ws.Range("E29").Value = "  -0.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "679.05"
$ws.Range("E30").Value = "  -4.08%  "
$ws.Range("E31").Value = "  -5.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.79"
$ws.Range("E32").Value = "  +12.33%  "
$ws.Range("E33").Value = "  -3.61%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "66.82"
$ws.Range("E34").Value = "  +8.72%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.445"
$ws.Range("E35").Value = "  -5.84%  "
$ws.Range("E36").Value = "  -8.75%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "39.38"
$ws.Range("E37").Value = "  -3.74%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.46"
$ws.Range("E38").Value = "  +12.64%  "
$ws.Range("E39").Value = "  -1.53%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  -0.09%  "
$ws.Range("E41").Value = "  -0.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.16"
$ws.Range("E42").Value = "  +5.02%  "
$ws.Range("E43").Value = "  -3.57%  "
$ws.Range("E44").Value = "  -4.24%  "
$ws.Range("E45").Value = "  +1.83%  "
$ws.Range("E46").Value = "  -1.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.99"
$ws.Range("E47").Value = "  -2.81%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.000272"
$ws.Range("E48").Value = "  +12.15%  "
$ws.Range("B49").Value = "LidoDAOToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.26"
$ws.Range("E49").Value = "  -3.64%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "143.27"
$ws.Range("E50").Value = "  -0.46%  "
$ws.Range("D51").Value = "0.0₆0337"
$ws.Range("E51").Value = "  -7.77%  "
